# Switched to EPPlus instead of Excel.Interop
#
# - Rename the "Process" sheet to "burp" (Excel auto-updates the
#   formula references on rename).
# - Input!D4: 8 -> 9
# - Output!C5 formula: wrap the VLOOKUP/multiplication in extra
#   parentheses and force an exact-match lookup (4th arg FALSE).
# - Restore the cell-cursor positions that the recorded macro captured
#   on each sheet (Input -> E4, burp -> B5), while leaving Output as
#   the active sheet/tab, matching the original workbook.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("Input")
$wsProcess = $wb.Worksheets.Item("Process")
$wsOutput = $wb.Worksheets.Item("Output")

# Rename "Process" -> "burp"
$wsProcess.Name = "burp"
$wsBurp = $wsProcess

# Input sheet: update value + selection
$wsInput.Activate()
$wsInput.Range("D4").Value = 9
$wsInput.Range("E4").Select()

# burp sheet: update selection only
$wsBurp.Activate()
$wsBurp.Range("B5").Select()

# Output sheet: update formula, then re-activate so it stays the
# active tab (matches the original file's state)
$wsOutput.Activate()
$wsOutput.Range("C5").Formula = "=((VLOOKUP(Input!D5,burp!A2:B5,2,FALSE))*C3)/1000"
$wsOutput.Range("C5").Select()
